$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- sheet1 ---
$ws1.Range("A2").Value = "Última actualización: 06:59:44"
$ws1.Range("A3").Value = "Total filas: 59"
$ws1.Range("A37").Value = "06:59:44"
$ws1.Range("C37").Value = "16_SANTA ANA"
$ws1.Range("D37").Value = 17
$ws1.Range("B38").Value = "07:16"
$ws1.Range("C38").Value = "11_ETCHEVERRY"
$ws1.Range("D38").Value = 32
$ws1.Range("A39").Value = "06:44:15"
$ws1.Range("B39").Value = "07:17"
$ws1.Range("C39").Value = "16_SANTA ANA"
$ws1.Range("D39").Value = 33
$ws1.Range("B40").Value = "07:21"
$ws1.Range("C40").Value = "26_HERNANDEZ"
$ws1.Range("D40").Value = 110
$ws1.Range("A41").Value = "05:31:23"
$ws1.Range("B41").Value = "07:23"
$ws1.Range("D41").Value = 112
$ws1.Range("A42").Value = "06:59:44"
$ws1.Range("B42").Value = "07:24"
$ws1.Range("C42").Value = "10_OLMOS"
$ws1.Range("D42").Value = 25
$ws1.Range("A43").Value = "06:44:15"
$ws1.Range("B43").Value = "07:25"
$ws1.Range("C43").Value = "10_OLMOS"
$ws1.Range("D43").Value = 41
$ws1.Range("B44").Value = "07:31"
$ws1.Range("C44").Value = "11_ETCHEVERRY"
$ws1.Range("D44").Value = 90
$ws1.Range("A45").Value = "06:01:37"
$ws1.Range("B45").Value = "07:31"
$ws1.Range("C45").Value = "16_SANTA ANA"
$ws1.Range("D45").Value = 90
$ws1.Range("A46").Value = "06:44:15"
$ws1.Range("B46").Value = "07:32"
$ws1.Range("C46").Value = "11_ETCHEVERRY"
$ws1.Range("D46").Value = 48
$ws1.Range("A47").Value = "06:01:37"
$ws1.Range("B47").Value = "07:32"
$ws1.Range("C47").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D47").Value = 91
$ws1.Range("B48").Value = "07:36"
$ws1.Range("C48").Value = "27_EL RETIRO"
$ws1.Range("D48").Value = 95
$ws1.Range("A49").Value = "06:44:15"
$ws1.Range("B49").Value = "07:37"
$ws1.Range("C49").Value = "27_EL RETIRO"
$ws1.Range("D49").Value = 53
$ws1.Range("A50").Value = "06:01:37"
$ws1.Range("B50").Value = "07:39"
$ws1.Range("C50").Value = "10_OLMOS"
$ws1.Range("D50").Value = 98
$ws1.Range("A51").Value = "06:59:44"
$ws1.Range("B51").Value = "07:45"
$ws1.Range("C51").Value = "11_ETCHEVERRY"
$ws1.Range("D51").Value = 46
$ws1.Range("A52").Value = "06:01:37"
$ws1.Range("B52").Value = "07:47"
$ws1.Range("C52").Value = "14_ABASTO"
$ws1.Range("D52").Value = 106
$ws1.Range("B53").Value = "07:48"
$ws1.Range("C53").Value = "14_ABASTO"
$ws1.Range("D53").Value = 64
$ws1.Range("A54").Value = "06:01:37"
$ws1.Range("B54").Value = "07:51"
$ws1.Range("C54").Value = "215D_EL PATO"
$ws1.Range("D54").Value = 110
$ws1.Range("B55").Value = "08:04"
$ws1.Range("C55").Value = "23_HERNANDEZ"
$ws1.Range("D55").Value = 80
$ws1.Range("A56").Value = "06:59:44"
$ws1.Range("B56").Value = "08:05"
$ws1.Range("C56").Value = "23_HERNANDEZ"
$ws1.Range("D56").Value = 66
$ws1.Range("B57").Value = "08:12"
$ws1.Range("C57").Value = "15_ABASTO"
$ws1.Range("D57").Value = 88
$ws1.Range("B58").Value = "08:21"
$ws1.Range("C58").Value = "26_HERNANDEZ"
$ws1.Range("D58").Value = 97
$ws1.Range("A59").Value = "06:59:44"
$ws1.Range("B59").Value = "08:22"
$ws1.Range("C59").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D59").Value = 83
$ws1.Range("E59").Value = "LP1912"
$ws1.Range("A60").Value = "06:44:15"
$ws1.Range("B60").Value = "08:23"
$ws1.Range("C60").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D60").Value = 99
$ws1.Range("E60").Value = "LP1912"
$ws1.Range("A61").Value = "06:44:15"
$ws1.Range("B61").Value = "08:23"
$ws1.Range("C61").Value = "215B_EL PATO"
$ws1.Range("D61").Value = 99
$ws1.Range("E61").Value = "LP1912"
$ws1.Range("A62").Value = "06:44:15"
$ws1.Range("B62").Value = "08:27"
$ws1.Range("C62").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D62").Value = 103
$ws1.Range("E62").Value = "LP1912"
$ws1.Range("A63").Value = "06:44:15"
$ws1.Range("B63").Value = "08:42"
$ws1.Range("C63").Value = "81_EL PELIGRO"
$ws1.Range("D63").Value = 118
$ws1.Range("E63").Value = "LP1912"
$ws1.Range("A64").Value = "06:59:44"
$ws1.Range("B64").Value = "08:54"
$ws1.Range("C64").Value = "17_ROMERO"
$ws1.Range("D64").Value = 115
$ws1.Range("E64").Value = "LP1912"

# --- sheet2 ---
$ws2.Range("A2").Value = "Última actualización: 06:59:44"

# --- sheet3 ---
$ws3.Range("A2").Value = "Última actualización: 06:59:44"
$ws3.Range("A3").Value = "Total filas: 14"
$ws3.Range("A15").Value = "06:59:44"
$ws3.Range("B15").Value = "07:39"
$ws3.Range("C15").Value = "215A_LA PLATA"
$ws3.Range("D15").Value = 40
$ws3.Range("E15").Value = "L6173"
$ws3.Range("B16").Value = "08:07"
$ws3.Range("C16").Value = "215C_LA PLATA"
$ws3.Range("D16").Value = 83
$ws3.Range("E16").Value = "L6203"
$ws3.Range("A17").Value = "06:59:44"
$ws3.Range("B17").Value = "08:09"
$ws3.Range("C17").Value = "215C_LA PLATA"
$ws3.Range("D17").Value = 70
$ws3.Range("E17").Value = "L6203"
$ws3.Range("A18").Value = "06:44:15"
$ws3.Range("B18").Value = "08:31"
$ws3.Range("C18").Value = "215A_LA PLATA"
$ws3.Range("D18").Value = 107
$ws3.Range("E18").Value = "L6173"
$ws3.Range("A19").Value = "06:59:44"
$ws3.Range("B19").Value = "08:35"
$ws3.Range("C19").Value = "215A_LA PLATA"
$ws3.Range("D19").Value = 96
$ws3.Range("E19").Value = "L6173"
